$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values below restore the D (Price) and E (Volume 1h) columns, and in two cases
# (rows 13-14) swap Chainlink/Avalanche back to their updated ranking order,
# matching the latest cryptos snapshot pulled by the scheduled GitHub Action.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.057.44"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.14"

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.98"
$ws.Range("E5").Value = "  +0.07%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.96"
$ws.Range("E7").Value = "  +0.41%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("E10").Value = "  +2.98%  "

# Row 11
$ws.Range("E11").Value = "  +0.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.400.19"
$ws.Range("E12").Value = "  +2.50%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.30"
$ws.Range("E13").Value = "  +3.96%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.62"
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("E15").Value = "  +5.92%  "

# Row 16
$ws.Range("E16").Value = "  +0.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.083.92"
$ws.Range("E17").Value = "  +2.88%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.005.57"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("E19").Value = "  +0.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.06"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.93"
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$ws.Range("E23").Value = "  +0.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("E25").Value = "  +3.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.94"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("E28").Value = "  +1.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +4.56%  "

# Row 31
$ws.Range("E31").Value = "  -0.48%  "

# Row 32
$ws.Range("E32").Value = "  +10.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.68"
$ws.Range("E33").Value = "  +2.99%  "

# Row 34
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("E35").Value = "  -0.45%  "

# Row 36
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("E37").Value = "  +4.75%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  +7.55%  "

# Row 39
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.19"
$ws.Range("E40").Value = "  +4.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.544.74"
$ws.Range("E41").Value = "  +1.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.87"
$ws.Range("E42").Value = "  +3.39%  "

# Row 43
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0903"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("E46").Value = "  +4.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -0.03%  "

# Row 48
$ws.Range("E48").Value = "  +1.01%  "

# Row 49
$ws.Range("E49").Value = "  +1.52%  "

# Row 50
$ws.Range("E50").Value = "  +0.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.286.46"
$ws.Range("E51").Value = "  +2.51%  "
